# "Phụ cấp" (allowance) is only paid at LONG XUYÊN. Remove the stray
# "Phụ cấp tại CẦN THƠ" and "Phụ cấp tại SÓC TRĂNG" rows from the salary
# breakdown sheet; every following row slides up to close the gap.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Remove "Phụ cấp tại CẦN THƠ" (row 3). Everything below shifts up one row.
$ws.Rows.Item(3).Delete()

# After the shift above, "Phụ cấp tại SÓC TRĂNG" is now row 23. Remove it too.
$ws.Rows.Item(23).Delete()

# The two rows whose labels land on the old "Lương cơ bản tại ..." slots
# (now rows 4 and 24) keep their originally-blank amount cell.
$ws.Range("B4").ClearContents()
$ws.Range("B24").ClearContents()
